$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Populate the "when to consider" (column D) helper values that were left
# blank for most feature rows. Categorical / binary features are marked 2,
# the two purely numeric features (age, absences) are marked 1. Rows that
# already carried a "D" value (2, 33, 34) are left untouched.
$twoRows = @(3, 4, 6, 7, 8, 11, 12, 13, 14, 18, 19, 20, 21, 22, 23, 24, 25)
foreach ($r in $twoRows) {
    $ws.Cells.Item($r, 4).Value = 2
}

$oneRows = @(5, 32)
foreach ($r in $oneRows) {
    $ws.Cells.Item($r, 4).Value = 1
}

# Append a small "Notes:" section below the table.
$ws.Range("B38").Value = "Notes:"
$ws.Range("B39").Value = "gender to be encoded in 0 and 1"

# Match the author's final selection/view state.
$ws.Range("F11").Select()
